$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update 想去人数 (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 386
$wsExhibit.Range("F3").Value = 2164

# Sheet "全部类型" (All types): same two rows appear again, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 386
$wsAll.Range("F7").Value = 2164
